$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 65.13500000000001
$ws.Range("D2").Value = 65.13500000000001
$ws.Range("E2").Value = 2.65552286
$ws.Range("F2").Value = 0.00638267
$ws.Range("G2").Value = 0.40896793
$ws.Range("H2").Value = 26.81933161
$ws.Range("I2").Value = 9.882126867806921
$ws.Range("J2").Value = 9.882126867806921
$ws.Range("K2").Value = 0.4098975460207761
$ws.Range("L2").Value = 0.001314625197731471
$ws.Range("M2").Value = 0.0699948853230469
$ws.Range("N2").Value = 6.839346901840893
$ws.Range("C3").Value = 111.549
$ws.Range("D3").Value = 111.549
$ws.Range("E3").Value = 1.56725154
$ws.Range("F3").Value = 0.0025236
$ws.Range("G3").Value = 0.27170185
$ws.Range("H3").Value = 30.24882202
$ws.Range("I3").Value = 21.52526273790719
$ws.Range("J3").Value = 21.52526273790719
$ws.Range("K3").Value = 0.2823447014375482
$ws.Range("L3").Value = 0.0008076659980241612
$ws.Range("M3").Value = 0.07151565521491303
$ws.Range("N3").Value = 9.985612663913578
$ws.Range("C4").Value = 42.358
$ws.Range("D4").Value = 84.68600000000001
$ws.Range("E4").Value = 2.10920756
$ws.Range("F4").Value = 0.01012117
$ws.Range("G4").Value = 0.20770927
$ws.Range("H4").Value = 8.987111659999998
$ws.Range("I4").Value = 9.661220170101455
$ws.Range("J4").Value = 19.31570605751963
$ws.Range("K4").Value = 0.5210583274160543
$ws.Range("L4").Value = 0.00217636336317918
$ws.Range("M4").Value = 0.04012654648109652
$ws.Range("N4").Value = 3.256796629911522
$ws.Range("C5").Value = 64.194
$ws.Range("D5").Value = 125.013
$ws.Range("E5").Value = 1.39997725
$ws.Range("F5").Value = 0.00471593
$ws.Range("G5").Value = 0.14674244
$ws.Range("H5").Value = 9.49293346
$ws.Range("I5").Value = 13.00813497683567
$ws.Range("J5").Value = 23.22611210003265
$ws.Range("K5").Value = 0.2687227255614431
$ws.Range("L5").Value = 0.001483955882147498
$ws.Range("M5").Value = 0.04011667138197247
$ws.Range("N5").Value = 3.44288405865803
$ws.Range("C6").Value = 25.737
$ws.Range("D6").Value = 102.866
$ws.Range("E6").Value = 1.77775795
$ws.Range("F6").Value = 0.01591173
$ws.Range("G6").Value = 0.09935664
$ws.Range("H6").Value = 2.67377001
$ws.Range("I6").Value = 7.126894388415169
$ws.Range("J6").Value = 28.49712494498577
$ws.Range("K6").Value = 0.5219397470318093
$ws.Range("L6").Value = 0.004499511140958271
$ws.Range("M6").Value = 0.03063395816252126
$ws.Range("N6").Value = 1.376432091603249
$ws.Range("C7").Value = 34.39
$ws.Range("D7").Value = 126.456
$ws.Range("E7").Value = 1.38925213
$ws.Range("F7").Value = 0.008069079999999999
$ws.Range("G7").Value = 0.06744921999999999
$ws.Range("H7").Value = 2.36027697
$ws.Range("I7").Value = 7.406669883203838
$ws.Range("J7").Value = 24.60344271425024
$ws.Range("K7").Value = 0.2816911078532773
$ws.Range("L7").Value = 0.002212040523335158
$ws.Range("M7").Value = 0.01694856043717261
$ws.Range("N7").Value = 0.9204318390499407
$ws.Range("C8").Value = 17.576
$ws.Range("D8").Value = 105.307
$ws.Range("E8").Value = 1.79116768
$ws.Range("F8").Value = 0.02023423
$ws.Range("G8").Value = 0.05714994
$ws.Range("H8").Value = 1.07914325
$ws.Range("I8").Value = 5.845031085322677
$ws.Range("J8").Value = 35.01292126709303
$ws.Range("K8").Value = 0.6156535921317461
$ws.Range("L8").Value = 0.006628597703482586
$ws.Range("M8").Value = 0.02103901701680555
$ws.Range("N8").Value = 0.7032740785929468
$ws.Range("C9").Value = 22.651
$ws.Range("D9").Value = 116.014
$ws.Range("E9").Value = 1.5195438
$ws.Range("F9").Value = 0.01131364
$ws.Range("G9").Value = 0.04170462
$ws.Range("H9").Value = 0.9718632
$ws.Range("I9").Value = 5.202717296496829
$ws.Range("J9").Value = 23.9911103823349
$ws.Range("K9").Value = 0.3162163444847202
$ws.Range("L9").Value = 0.00348032530298744
$ws.Range("M9").Value = 0.01301746547769815
$ws.Range("N9").Value = 0.4624244537393932
$ws.Range("C10").Value = 13.279
$ws.Range("D10").Value = 106.027
$ws.Range("E10").Value = 1.80481392
$ws.Range("F10").Value = 0.02214199
$ws.Range("G10").Value = 0.0354577
$ws.Range("H10").Value = 0.51102201
$ws.Range("I10").Value = 4.611336403194136
$ws.Range("J10").Value = 36.83388393865091
$ws.Range("K10").Value = 0.6766637493802656
$ws.Range("L10").Value = 0.007590792769847875
$ws.Range("M10").Value = 0.01413135320673337
$ws.Range("N10").Value = 0.3517711518190655
$ws.Range("C11").Value = 16.439
$ws.Range("D11").Value = 102.816
$ws.Range("E11").Value = 1.72266275
$ws.Range("F11").Value = 0.01285501
$ws.Range("G11").Value = 0.02577802
$ws.Range("H11").Value = 0.4382587099999999
$ws.Range("I11").Value = 4.040371716018663
$ws.Range("J11").Value = 22.76608087568957
$ws.Range("K11").Value = 0.3718583440531885
$ws.Range("L11").Value = 0.003959105924273764
$ws.Range("M11").Value = 0.008166565353272407
$ws.Range("N11").Value = 0.2145351145855292
$ws.Range("C12").Value = 10.123
$ws.Range("D12").Value = 101.044
$ws.Range("E12").Value = 1.91592889
$ws.Range("F12").Value = 0.02498578
$ws.Range("G12").Value = 0.02418259
$ws.Range("H12").Value = 0.26764387
$ws.Range("I12").Value = 3.756325462600516
$ws.Range("J12").Value = 37.58021203155108
$ws.Range("K12").Value = 0.7386340800974114
$ws.Range("L12").Value = 0.008863485746091724
$ws.Range("M12").Value = 0.009957136453732982
$ws.Range("N12").Value = 0.2000475682708845
$ws.Range("C13").Value = 12.954
$ws.Range("D13").Value = 91.63500000000001
$ws.Range("E13").Value = 1.93368706
$ws.Range("F13").Value = 0.0143713
$ws.Range("G13").Value = 0.0182102
$ws.Range("H13").Value = 0.24611213
$ws.Range("I13").Value = 3.3608301324493
$ws.Range("J13").Value = 20.58882402359772
$ws.Range("K13").Value = 0.4181174380136848
$ws.Range("L13").Value = 0.004688753984154417
$ws.Range("M13").Value = 0.006375824958693887
$ws.Range("N13").Value = 0.1400917751533927
